# Apply the "pi 09. 04. 2021" daily update to the Slovakia Covid DailyStats sheet.
# Updates the AgTests (F) and AgPosit (G) columns for rows 362-399
# (dates 2021-03-01 through 2021-04-07) to reflect the refreshed open-data export.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F362").Value = 228355
$ws.Range("G362").Value = 3182
$ws.Range("F363").Value = 187818
$ws.Range("G363").Value = 2759
$ws.Range("F364").Value = 167716
$ws.Range("G364").Value = 2470
$ws.Range("F365").Value = 183780
$ws.Range("G365").Value = 2386
$ws.Range("F366").Value = 340044
$ws.Range("G366").Value = 2846
$ws.Range("F367").Value = 765921
$ws.Range("G367").Value = 3922
$ws.Range("F369").Value = 233518
$ws.Range("G369").Value = 2593
$ws.Range("F370").Value = 182217
$ws.Range("G370").Value = 2043
$ws.Range("F371").Value = 159310
$ws.Range("G371").Value = 1954
$ws.Range("F372").Value = 178757
$ws.Range("G372").Value = 1856
$ws.Range("F373").Value = 348167
$ws.Range("G373").Value = 2371
$ws.Range("F374").Value = 771238
$ws.Range("G374").Value = 3419
$ws.Range("F376").Value = 221053
$ws.Range("G376").Value = 2225
$ws.Range("F377").Value = 176750
$ws.Range("G377").Value = 1822
$ws.Range("F378").Value = 157243
$ws.Range("G378").Value = 1545
$ws.Range("F379").Value = 179021
$ws.Range("F380").Value = 343690
$ws.Range("F381").Value = 743351
$ws.Range("G381").Value = 2685
$ws.Range("F382").Value = 357556
$ws.Range("F383").Value = 220610
$ws.Range("G383").Value = 1759
$ws.Range("F384").Value = 171828
$ws.Range("G384").Value = 1511
$ws.Range("F385").Value = 150713
$ws.Range("G385").Value = 1403
$ws.Range("F386").Value = 182343
$ws.Range("G386").Value = 1364
$ws.Range("F387").Value = 350636
$ws.Range("G387").Value = 1665
$ws.Range("F388").Value = 727907
$ws.Range("G388").Value = 2197
$ws.Range("F389").Value = 353443
$ws.Range("G389").Value = 1304
$ws.Range("F390").Value = 219669
$ws.Range("G390").Value = 1508
$ws.Range("F391").Value = 176676
$ws.Range("G391").Value = 1214
$ws.Range("F392").Value = 217395
$ws.Range("G392").Value = 1197
$ws.Range("F393").Value = 294377
$ws.Range("G393").Value = 1179
$ws.Range("F394").Value = 161015
$ws.Range("G394").Value = 613
$ws.Range("F395").Value = 731930
$ws.Range("G395").Value = 1903
$ws.Range("F396").Value = 163208
$ws.Range("G396").Value = 546
$ws.Range("F397").Value = 104418
$ws.Range("G397").Value = 617
$ws.Range("F398").Value = 285706
$ws.Range("G398").Value = 1419
$ws.Range("F399").Value = 190457
$ws.Range("G399").Value = 942
